$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 69: the trailing empty "inlineStr" cell in column G is removed entirely
# (Dismissed charge with no finding recorded in G).
$ws.Range("G69").ClearContents()

# Charge-grid rows moved in from elsewhere: same case (21TRD09437 / Bunner)
# re-entered as "No Contest" pleas resulting in "Guilty" findings.
$newRows = @(
    @{ Row = 70; A = "21TRD09437"; B = "Bunner"; C = "DUS"; D = "4510.11"; E = "M1"; F = "No Contest"; G = "Guilty"; H = 0; I = "0" },
    @{ Row = 71; A = "21TRD09437"; B = "Bunner"; C = "1ST SPEED 1 YR SCHOOL >35MPHM4"; D = "4511.21B1A"; E = "M4"; F = "No Contest"; G = "Guilty"; H = 0; I = "0" },
    @{ Row = 72; A = "21TRD09437"; B = "Bunner"; C = "RECKLESS OPERATION 1ST IN 1 YR"; D = "4511.20"; E = "MM"; F = "No Contest"; G = "Guilty"; H = 0; I = "0" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C

    # Column D holds statute codes that look numeric (e.g. "4511.20") but
    # must stay TEXT so formatting/trailing zeros/precision survive.
    $ws.Range("D$row").NumberFormat = "@"
    $ws.Range("D$row").Value = $r.D

    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
    $ws.Range("H$row").Value = $r.H

    # Column I holds the numeric-looking value "0" as TEXT (not a number),
    # matching the rest of the sheet's inline-string convention for this column.
    $ws.Range("I$row").NumberFormat = "@"
    $ws.Range("I$row").Value = $r.I
}
